# Updates the "cryptos" price list on Sheet1 (rows 2-51) to the latest
# scrape: refreshed prices/volume deltas in columns D/E for every coin,
# plus two rank swaps (Litecoin <-> ShibaInu at rows 19/20, and
# Aave <-> Gas at rows 44/45) where name/link/price/volume move together.
#
# Column D holds price text that sometimes looks like a plain decimal
# (e.g. "55.42", "237.00"). Excel's COM layer auto-converts such literals
# to numbers on assignment (e.g. "237.00" -> 237), which would silently
# strip the significant trailing zero / change the cell's stored type
# from Text to Number. Setting NumberFormat to "@" (Text) before writing
# the value keeps these cells literal text, matching the source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


# Row 2
$ws.Cells.Item(2, 4).NumberFormat = '@'
$ws.Cells.Item(2, 4).Value = '37.136.36'
$ws.Cells.Item(2, 5).Value = '  +1.38%  '

# Row 3
$ws.Cells.Item(3, 4).NumberFormat = '@'
$ws.Cells.Item(3, 4).Value = '2.058.32'
$ws.Cells.Item(3, 5).Value = '  -2.63%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.03%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '248.75'

# Row 6
$ws.Cells.Item(6, 5).Value = '  -1.23%  '

# Row 7
$ws.Cells.Item(7, 5).Value = '  -0.04%  '

# Row 8
$ws.Cells.Item(8, 4).NumberFormat = '@'
$ws.Cells.Item(8, 4).Value = '55.42'
$ws.Cells.Item(8, 5).Value = '  +17.15%  '

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '61.47'
$ws.Cells.Item(9, 5).Value = '  +1.38%  '

# Row 10
$ws.Cells.Item(10, 5).Value = '  +1.85%  '

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '0.0795'
$ws.Cells.Item(11, 5).Value = '  +6.47%  '

# Row 12
$ws.Cells.Item(12, 5).Value = '  +5.50%  '

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = '@'
$ws.Cells.Item(13, 4).Value = '15.19'
$ws.Cells.Item(13, 5).Value = '  +6.44%  '

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '2.354.48'
$ws.Cells.Item(14, 5).Value = '  -2.69%  '

# Row 15
$ws.Cells.Item(15, 5).Value = '  -1.67%  '

# Row 16
$ws.Cells.Item(16, 5).Value = '  +2.46%  '

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '2.061.16'
$ws.Cells.Item(17, 5).Value = '  -2.42%  '

# Row 18
$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '37.087.47'
$ws.Cells.Item(18, 5).Value = '  +1.28%  '

# Row 19
$ws.Cells.Item(19, 2).Value = 'ShibaInu'
$ws.Cells.Item(19, 3).Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '0.0₃0908'
$ws.Cells.Item(19, 5).Value = '  +8.71%  '

# Row 20
$ws.Cells.Item(20, 2).Value = 'Litecoin'
$ws.Cells.Item(20, 3).Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '72.35'
$ws.Cells.Item(20, 5).Value = '  -1.53%  '

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '14.19'
$ws.Cells.Item(21, 5).Value = '  +7.32%  '

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '5.34'
$ws.Cells.Item(22, 5).Value = '  +2.84%  '

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '237.00'
$ws.Cells.Item(23, 5).Value = '  -1.46%  '

# Row 24
$ws.Cells.Item(24, 5).Value = '  -0.01%  '

# Row 25
$ws.Cells.Item(25, 5).Value = '  -2.28%  '

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = '@'
$ws.Cells.Item(26, 4).Value = '169.93'
$ws.Cells.Item(26, 5).Value = '  -1.41%  '

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '9.03'
$ws.Cells.Item(27, 5).Value = '  -1.75%  '

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '20.16'
$ws.Cells.Item(28, 5).Value = '  -6.73%  '

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '1.98'
$ws.Cells.Item(29, 5).Value = '  -1.69%  '

# Row 30
$ws.Cells.Item(30, 5).Value = '  -0.12%  '

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '4.54'
$ws.Cells.Item(31, 5).Value = '  +1.05%  '

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '1.06'
$ws.Cells.Item(32, 5).Value = '  +10.65%  '

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '0.0625'
$ws.Cells.Item(33, 5).Value = '  +4.12%  '

# Row 34
$ws.Cells.Item(34, 5).Value = '  +5.15%  '

# Row 35
$ws.Cells.Item(35, 5).Value = '  +0.04%  '

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '0.0869'
$ws.Cells.Item(36, 5).Value = '  -7.16%  '

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '2.27'
$ws.Cells.Item(37, 5).Value = '  -3.73%  '

# Row 38
$ws.Cells.Item(38, 5).Value = '  -5.75%  '

# Row 39
$ws.Cells.Item(39, 5).Value = '  +1.24%  '

# Row 40
$ws.Cells.Item(40, 5).Value = '  +21.86%  '

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '18.14'
$ws.Cells.Item(41, 5).Value = '  +13.74%  '

# Row 42
$ws.Cells.Item(42, 5).Value = '  -0.78%  '

# Row 44
$ws.Cells.Item(44, 2).Value = 'Gas'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '14.86'
$ws.Cells.Item(44, 5).Value = '  -50.10%  '

# Row 45
$ws.Cells.Item(45, 2).Value = 'Aave'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '96.22'
$ws.Cells.Item(45, 5).Value = '  -2.45%  '

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '4.33'
$ws.Cells.Item(46, 5).Value = '  +48.49%  '

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '2.79'
$ws.Cells.Item(47, 5).Value = '  +0.14%  '

# Row 48
$ws.Cells.Item(48, 5).Value = '  +5.71%  '

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '1.298.11'
$ws.Cells.Item(49, 5).Value = '  -3.79%  '

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '2.91'
$ws.Cells.Item(50, 5).Value = '  +2.75%  '

# Row 51
$ws.Cells.Item(51, 5).Value = '  -5.32%  '
